$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the value in B6 (previously "厉鬼") while keeping its style
$ws.Range("B6").ClearContents()

# Update the selection to B6 (was F6)
$ws.Range("B6").Select()
